$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(4, -1, -1, 0, 0, -1),
    @(-1, 4, -1, -1, 0, 2),
    @(-1, -1, 4, -1, -1, 6),
    @(0, -1, -1, 4, -1, 2),
    @(0, 0, -1, -1, 4, 1)
)

for ($r = 0; $r -lt 5; $r++) {
    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $values[$r][$c]
    }
}

$ws.Range("F8").Select()
